$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-05 Sunday", 2) | Out-Null
$d.Content.Find.Execute("802÷2=401, 0", $true, $false, $false, $false, $false, $true, 1, $false, "812÷4=203, 0", 2) | Out-Null
$d.Content.Find.Execute("127÷8=15, 7", $true, $false, $false, $false, $false, $true, 1, $false, "681÷9=75, 6", 2) | Out-Null
$d.Content.Find.Execute("964÷3=321, 1", $true, $false, $false, $false, $false, $true, 1, $false, "123÷5=24, 3", 2) | Out-Null
$d.Content.Find.Execute("748÷4=187, 0", $true, $false, $false, $false, $false, $true, 1, $false, "896÷6=149, 2", 2) | Out-Null
$d.Content.Find.Execute("355÷9=39, 4", $true, $false, $false, $false, $false, $true, 1, $false, "156÷6=26, 0", 2) | Out-Null
$d.Content.Find.Execute("281÷7=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "716÷5=143, 1", 2) | Out-Null
$d.Content.Find.Execute("370÷8=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "217÷9=24, 1", 2) | Out-Null
$d.Content.Find.Execute("851÷5=170, 1", $true, $false, $false, $false, $false, $true, 1, $false, "285÷3=95, 0", 2) | Out-Null
$d.Content.Find.Execute("640÷7=91, 3", $true, $false, $false, $false, $false, $true, 1, $false, "841÷2=420, 1", 2) | Out-Null
$d.Content.Find.Execute("435÷8=54, 3", $true, $false, $false, $false, $false, $true, 1, $false, "437÷4=109, 1", 2) | Out-Null
$d.Content.Find.Execute("210÷9=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "418÷5=83, 3", 2) | Out-Null
$d.Content.Find.Execute("298÷5=59, 3", $true, $false, $false, $false, $false, $true, 1, $false, "533÷9=59, 2", 2) | Out-Null
$d.Content.Find.Execute("247÷6=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "832÷9=92, 4", 2) | Out-Null
$d.Content.Find.Execute("535÷7=76, 3", $true, $false, $false, $false, $false, $true, 1, $false, "832÷4=208, 0", 2) | Out-Null
$d.Content.Find.Execute("129÷7=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "341÷8=42, 5", 2) | Out-Null
$d.Content.Find.Execute("217÷2=108, 1", $true, $false, $false, $false, $false, $true, 1, $false, "198÷3=66, 0", 2) | Out-Null
$d.Content.Find.Execute("942÷7=134, 4", $true, $false, $false, $false, $false, $true, 1, $false, "410÷2=205, 0", 2) | Out-Null
$d.Content.Find.Execute("578÷8=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "288÷9=32, 0", 2) | Out-Null
$d.Content.Find.Execute("555÷9=61, 6", $true, $false, $false, $false, $false, $true, 1, $false, "281÷6=46, 5", 2) | Out-Null
$d.Content.Find.Execute("828÷6=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "484÷8=60, 4", 2) | Out-Null
$d.Content.Find.Execute("711÷9=79, 0", $true, $false, $false, $false, $false, $true, 1, $false, "847÷4=211, 3", 2) | Out-Null
$d.Content.Find.Execute("244÷9=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "941÷6=156, 5", 2) | Out-Null
$d.Content.Find.Execute("303÷6=50, 3", $true, $false, $false, $false, $false, $true, 1, $false, "725÷5=145, 0", 2) | Out-Null
$d.Content.Find.Execute("768÷6=128, 0", $true, $false, $false, $false, $false, $true, 1, $false, "740÷9=82, 2", 2) | Out-Null
$d.Content.Find.Execute("969÷2=484, 1", $true, $false, $false, $false, $false, $true, 1, $false, "736÷6=122, 4", 2) | Out-Null
